$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new timesheet entry ------------------------------------------
# Date: 2023-01-19
$ws.Range("A10").Value = 44945
$ws.Range("A10").NumberFormat = "d-mmm"

# Time started: 8:43 AM
$ws.Range("B10").Value = 0.36319444444444443
$ws.Range("B10").NumberFormat = "h:mm"

# Time finished: 10:13 AM
$ws.Range("C10").Value = 0.42569444444444443
$ws.Range("C10").NumberFormat = "h:mm"

# Description
$ws.Range("D10").Value = "Restructured project. More backend for adding tutors. Worked on referential integrity of database."

# Total hours for the entry
$ws.Range("E10").Value = 1.5

# --- Row heights now that D9/D10 hold wrapped multi-line text -------------
$ws.Rows.Item(9).RowHeight = 28
$ws.Rows.Item(10).RowHeight = 42

# --- Selection / active cell, as last left by the author ------------------
$ws.Range("D12").Select()
